# Update NATMI edge-weight outputs for the Slitrk6-Ptprs sheet with refreshed TPM values.
# Ligand average expression for MuSCs (G12:G16) and receptor average expression for the
# ECs / Inflammatory-Mac / MuSCs / Resolving-Mac target clusters (M column) were recomputed
# upstream; all dependent specificity / edge-weight columns below are their refreshed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "I2" = 0.04018528407978143
    "J2" = 0.04018528407978143
    "M2" = 4.925988333333333
    "N2" = 14.777965
    "O2" = 0.05656988822582037
    "P2" = 0.05656988822582035
    "Q2" = 0.2466228899005555
    "R2" = 2.219606009105
    "S2" = 0.002273277028716074
    "T2" = 0.002273277028716074
    "I3" = 0.04018528407978143
    "J3" = 0.04018528407978143
    "O3" = 0.5464678959362861
    "P3" = 0.5464678959362861
    "S3" = 0.02195996763868009
    "T3" = 0.02195996763868009
    "I4" = 0.04018528407978143
    "J4" = 0.04018528407978143
    "M4" = 13.80191933333334
    "N4" = 41.40575800000001
    "O4" = 0.1585007882996995
    "P4" = 0.1585007882996994
    "Q4" = 0.691002292702889
    "R4" = 6.219020634326001
    "S4" = 0.006369399204692719
    "T4" = 0.006369399204692718
    "I5" = 0.04018528407978143
    "J5" = 0.04018528407978143
    "M5" = 12.18796133333333
    "N5" = 36.563884
    "O5" = 0.1399661476381804
    "P5" = 0.1399661476381803
    "Q5" = 0.6101984094608889
    "R5" = 5.491785685148
    "S5" = 0.005624579404392906
    "T5" = 0.005624579404392905
    "I6" = 0.04018528407978143
    "J6" = 0.04018528407978143
    "M6" = 8.576764333333333
    "N6" = 25.730293
    "O6" = 0.09849527990001386
    "P6" = 0.09849527990001385
    "Q6" = 0.4294014241912222
    "R6" = 3.864612817721
    "S6" = 0.003958060803299642
    "T6" = 0.003958060803299642
    "I7" = 0.8039569115253269
    "J7" = 0.8039569115253269
    "M7" = 4.925988333333333
    "N7" = 14.777965
    "O7" = 0.05656988822582037
    "P7" = 0.05656988822582035
    "Q7" = 4.933999632359444
    "R7" = 44.40599669123499
    "S7" = 0.0454797526233635
    "T7" = 0.04547975262336348
    "I8" = 0.8039569115253269
    "J8" = 0.8039569115253269
    "O8" = 0.5464678959362861
    "P8" = 0.5464678959362861
    "S8" = 0.4393366418646803
    "T8" = 0.4393366418646803
    "I9" = 0.8039569115253269
    "J9" = 0.8039569115253269
    "M9" = 13.80191933333334
    "N9" = 41.40575800000001
    "O9" = 0.1585007882996995
    "P9" = 0.1585007882996994
    "Q9" = 13.82436585480911
    "R9" = 124.419292693282
    "S9" = 0.1274278042357561
    "T9" = 0.127427804235756
    "I10" = 0.8039569115253269
    "J10" = 0.8039569115253269
    "M10" = 12.18796133333333
    "N10" = 36.563884
    "O10" = 0.1399661476381804
    "P10" = 0.1399661476381803
    "Q10" = 12.20778302111511
    "R10" = 109.870047190036
    "S10" = 0.1125267517732894
    "T10" = 0.1125267517732894
    "I11" = 0.8039569115253269
    "J11" = 0.8039569115253269
    "M11" = 8.576764333333333
    "N11" = 25.730293
    "O11" = 0.09849527990001386
    "P11" = 0.09849527990001385
    "Q11" = 8.590713011060778
    "R11" = 77.316417099547
    "S11" = 0.07918596102823774
    "T11" = 0.07918596102823774
    "G12" = 0.1941786666666666
    "H12" = 0.5825359999999999
    "I12" = 0.1558578043948917
    "J12" = 0.1558578043948917
    "M12" = 4.925988333333333
    "N12" = 14.777965
    "O12" = 0.05656988822582037
    "P12" = 0.05656988822582035
    "Q12" = 0.9565218465822219
    "R12" = 8.608696619239998
    "S12" = 0.008816858573740798
    "T12" = 0.008816858573740797
    "G13" = 0.1941786666666666
    "H13" = 0.5825359999999999
    "I13" = 0.1558578043948917
    "J13" = 0.1558578043948917
    "O13" = 0.5464678959362861
    "P13" = 0.5464678959362861
    "Q13" = 9.240047970967998
    "R13" = 83.16043173871199
    "S13" = 0.08517128643292571
    "T13" = 0.08517128643292572
    "G14" = 0.1941786666666666
    "H14" = 0.5825359999999999
    "I14" = 0.1558578043948917
    "J14" = 0.1558578043948917
    "M14" = 13.80191933333334
    "N14" = 41.40575800000001
    "O14" = 0.1585007882996995
    "P14" = 0.1585007882996994
    "Q14" = 2.680038293587556
    "R14" = 24.120344642288
    "S14" = 0.0247035848592507
    "T14" = 0.0247035848592507
    "G15" = 0.1941786666666666
    "H15" = 0.5825359999999999
    "I15" = 0.1558578043948917
    "J15" = 0.1558578043948917
    "M15" = 12.18796133333333
    "N15" = 36.563884
    "O15" = 0.1399661476381804
    "P15" = 0.1399661476381803
    "Q15" = 2.366642081091555
    "R15" = 21.299778729824
    "S15" = 0.02181481646049805
    "T15" = 0.02181481646049805
    "G16" = 0.1941786666666666
    "H16" = 0.5825359999999999
    "I16" = 0.1558578043948917
    "J16" = 0.1558578043948917
    "M16" = 8.576764333333333
    "N16" = 25.730293
    "O16" = 0.09849527990001386
    "P16" = 0.09849527990001385
    "Q16" = 1.665424662560889
    "R16" = 14.988821963048
    "S16" = 0.01535125806847647
    "T16" = 0.01535125806847647
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
